# Append the final row of sequential-memory-usage data to the sheet.
# (commit: "final output data and upgraded diagrams")
#
# The sheet currently holds data in rows 1 (headers) through 11
# (A1:L11). A new data row (row 12) is appended that only carries
# values in the last two columns, K and L, extending the used range
# to A1:L12. Excel/the workbook will automatically grow the sheet
# dimension and keep the existing charts (which reference the full
# column ranges) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K12").Value = 23919576
$ws.Range("L12").Value = 42370048
